$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 41793
$ws.Range("G2").Value = "暂时售罄"
$ws.Range("F5").Value = 9495
$ws.Range("F6").Value = 201
$ws.Range("F7").Value = 855
$ws.Range("F8").Value = 897
$ws.Range("F9").Value = 724
$ws.Range("F10").Value = 214
$ws.Range("F11").Value = 56
$ws.Range("F12").Value = 295
$ws.Range("F13").Value = 894
$ws.Range("F14").Value = 84
$ws.Range("F15").Value = 124
$ws.Range("F16").Value = 730
$ws.Range("F18").Value = 1402
$ws.Range("F20").Value = 658
$ws.Range("F21").Value = 699
$ws.Range("F22").Value = 457
$ws.Range("F24").Value = 731
$ws.Range("F26").Value = 248
$ws.Range("F27").Value = 62
$ws.Range("F28").Value = 498
$ws.Range("F29").Value = 521
$ws.Range("F31").Value = 239
$ws.Range("F32").Value = 925
$ws.Range("F33").Value = 14
$ws.Range("F35").Value = 94
$ws.Range("F36").Value = 213
$ws.Range("F37").Value = 144
$ws.Range("F38").Value = 391
$ws.Range("F39").Value = 1257
$ws.Range("F40").Value = 290
$ws.Range("F42").Value = 1235
$ws.Range("F43").Value = 375
$ws.Range("F45").Value = 13

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 214
$ws.Range("G3").Value = 178
$ws.Range("F5").Value = 4445
$ws.Range("F10").Value = 78
$ws.Range("F11").Value = 127
$ws.Range("F17").Value = 164
$ws.Range("F22").Value = 6

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2027
$ws.Range("F3").Value = 521
$ws.Range("F4").Value = 400

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2027
$ws.Range("F3").Value = 521
$ws.Range("F4").Value = 41793
$ws.Range("G4").Value = "暂时售罄"
$ws.Range("F7").Value = 214
$ws.Range("G7").Value = 178
$ws.Range("F10").Value = 9495
$ws.Range("F11").Value = 201
$ws.Range("F12").Value = 855
$ws.Range("F13").Value = 855
$ws.Range("F14").Value = 78
$ws.Range("F15").Value = 400
$ws.Range("F16").Value = 897
$ws.Range("F17").Value = 127
$ws.Range("F18").Value = 214
$ws.Range("F19").Value = 295
$ws.Range("F20").Value = 894
$ws.Range("F22").Value = 84
$ws.Range("F24").Value = 730
$ws.Range("F26").Value = 1402
$ws.Range("F27").Value = 658
$ws.Range("F28").Value = 699
$ws.Range("F29").Value = 457
$ws.Range("F31").Value = 732
$ws.Range("F33").Value = 62
$ws.Range("F34").Value = 498
$ws.Range("F36").Value = 239
$ws.Range("F37").Value = 925
$ws.Range("F38").Value = 14
$ws.Range("F40").Value = 94
$ws.Range("F41").Value = 213
$ws.Range("F42").Value = 391
$ws.Range("F43").Value = 1235
$ws.Range("F44").Value = 375
